$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.046494
$ws.Range("H2").Value = 0.139482
$ws.Range("I2").Value = 0.05587126560541624
$ws.Range("J2").Value = 0.05587126560541625
$ws.Range("M2").Value = 6.066605666666667
$ws.Range("N2").Value = 18.199817
$ws.Range("O2").Value = 0.8497846287916651
$ws.Range("P2").Value = 0.8497846287916652
$ws.Range("Q2").Value = 0.282060763866
$ws.Range("R2").Value = 2.538546874794
$ws.Range("S2").Value = 0.04747854270261916
$ws.Range("T2").Value = 0.04747854270261918

$ws.Range("G3").Value = 0.046494
$ws.Range("H3").Value = 0.139482
$ws.Range("I3").Value = 0.05587126560541624
$ws.Range("J3").Value = 0.05587126560541625
$ws.Range("O3").Value = 0.1196497582104962
$ws.Range("P3").Value = 0.1196497582104962
$ws.Range("Q3").Value = 0.039714182928
$ws.Range("R3").Value = 0.357427646352
$ws.Range("S3").Value = 0.006684983420602465
$ws.Range("T3").Value = 0.006684983420602466

$ws.Range("G4").Value = 0.046494
$ws.Range("H4").Value = 0.139482
$ws.Range("I4").Value = 0.05587126560541624
$ws.Range("J4").Value = 0.05587126560541625
$ws.Range("M4").Value = 0.1824346666666667
$ws.Range("N4").Value = 0.547304
$ws.Range("O4").Value = 0.02555468148257719
$ws.Range("P4").Value = 0.02555468148257719
$ws.Range("Q4").Value = 0.008482117391999999
$ws.Range("R4").Value = 0.07633905652799999
$ws.Range("S4").Value = 0.001427772396574882
$ws.Range("T4").Value = 0.001427772396574882

$ws.Range("G5").Value = 0.046494
$ws.Range("H5").Value = 0.139482
$ws.Range("I5").Value = 0.05587126560541624
$ws.Range("J5").Value = 0.05587126560541625
$ws.Range("M5").Value = 0.035773
$ws.Range("N5").Value = 0.107319
$ws.Range("O5").Value = 0.005010931515261538
$ws.Range("P5").Value = 0.005010931515261539
$ws.Range("Q5").Value = 0.001663229862
$ws.Range("R5").Value = 0.014969068758
$ws.Range("S5").Value = 0.0002799670856197283
$ws.Range("T5").Value = 0.0002799670856197283

$ws.Range("I6").Value = 0.7083039420562237
$ws.Range("J6").Value = 0.7083039420562237
$ws.Range("M6").Value = 6.066605666666667
$ws.Range("N6").Value = 18.199817
$ws.Range("O6").Value = 0.8497846287916651
$ws.Range("P6").Value = 0.8497846287916652
$ws.Range("Q6").Value = 3.575805000671222
$ws.Range("R6").Value = 32.182245006041
$ws.Range("S6").Value = 0.601905802471921
$ws.Range("T6").Value = 0.6019058024719212

$ws.Range("I7").Value = 0.7083039420562237
$ws.Range("J7").Value = 0.7083039420562237
$ws.Range("O7").Value = 0.1196497582104962
$ws.Range("P7").Value = 0.1196497582104962
$ws.Range("S7").Value = 0.08474839540656846
$ws.Range("T7").Value = 0.08474839540656848

$ws.Range("I8").Value = 0.7083039420562237
$ws.Range("J8").Value = 0.7083039420562237
$ws.Range("M8").Value = 0.1824346666666667
$ws.Range("N8").Value = 0.547304
$ws.Range("O8").Value = 0.02555468148257719
$ws.Range("P8").Value = 0.02555468148257719
$ws.Range("Q8").Value = 0.1075314317768889
$ws.Range("R8").Value = 0.967782885992
$ws.Range("S8").Value = 0.0181004816321006
$ws.Range("T8").Value = 0.01810048163210061

$ws.Range("I9").Value = 0.7083039420562237
$ws.Range("J9").Value = 0.7083039420562237
$ws.Range("M9").Value = 0.035773
$ws.Range("N9").Value = 0.107319
$ws.Range("O9").Value = 0.005010931515261538
$ws.Range("P9").Value = 0.005010931515261539
$ws.Range("Q9").Value = 0.02108547667633333
$ws.Range("R9").Value = 0.189769290087
$ws.Range("S9").Value = 0.003549262545633514
$ws.Range("T9").Value = 0.003549262545633514

$ws.Range("G10").Value = 0.1246316666666667
$ws.Range("H10").Value = 0.373895
$ws.Range("I10").Value = 0.1497683346491813
$ws.Range("J10").Value = 0.1497683346491813
$ws.Range("M10").Value = 6.066605666666667
$ws.Range("N10").Value = 18.199817
$ws.Range("O10").Value = 0.8497846287916651
$ws.Range("P10").Value = 0.8497846287916652
$ws.Range("Q10").Value = 0.756091175246111
$ws.Range("R10").Value = 6.804820577215
$ws.Range("S10").Value = 0.1272708286646004
$ws.Range("T10").Value = 0.1272708286646004

$ws.Range("G11").Value = 0.1246316666666667
$ws.Range("H11").Value = 0.373895
$ws.Range("I11").Value = 0.1497683346491813
$ws.Range("J11").Value = 0.1497683346491813
$ws.Range("O11").Value = 0.1196497582104962
$ws.Range("P11").Value = 0.1196497582104962
$ws.Range("Q11").Value = 0.1064577108577777
$ws.Range("R11").Value = 0.9581193977199999
$ws.Range("S11").Value = 0.01791974502836322
$ws.Range("T11").Value = 0.01791974502836322

$ws.Range("G12").Value = 0.1246316666666667
$ws.Range("H12").Value = 0.373895
$ws.Range("I12").Value = 0.1497683346491813
$ws.Range("J12").Value = 0.1497683346491813
$ws.Range("M12").Value = 0.1824346666666667
$ws.Range("N12").Value = 0.547304
$ws.Range("O12").Value = 0.02555468148257719
$ws.Range("P12").Value = 0.02555468148257719
$ws.Range("Q12").Value = 0.02273713656444444
$ws.Range("R12").Value = 0.20463422908
$ws.Range("S12").Value = 0.003827282088135856
$ws.Range("T12").Value = 0.003827282088135857

$ws.Range("G13").Value = 0.1246316666666667
$ws.Range("H13").Value = 0.373895
$ws.Range("I13").Value = 0.1497683346491813
$ws.Range("J13").Value = 0.1497683346491813
$ws.Range("M13").Value = 0.035773
$ws.Range("N13").Value = 0.107319
$ws.Range("O13").Value = 0.005010931515261538
$ws.Range("P13").Value = 0.005010931515261539
$ws.Range("Q13").Value = 0.004458448611666666
$ws.Range("R13").Value = 0.040126037505
$ws.Range("S13").Value = 0.0007504788680818189
$ws.Range("T13").Value = 0.0007504788680818193

$ws.Range("G14").Value = 0.05954566666666666
$ws.Range("H14").Value = 0.178637
$ws.Range("I14").Value = 0.07155529225243931
$ws.Range("J14").Value = 0.07155529225243933
$ws.Range("M14").Value = 6.066605666666667
$ws.Range("N14").Value = 18.199817
$ws.Range("O14").Value = 0.8497846287916651
$ws.Range("P14").Value = 0.8497846287916652
$ws.Range("Q14").Value = 0.3612400788254445
$ws.Range("R14").Value = 3.251160709429
$ws.Range("S14").Value = 0.06080658746481825
$ws.Range("T14").Value = 0.06080658746481827

$ws.Range("G15").Value = 0.05954566666666666
$ws.Range("H15").Value = 0.178637
$ws.Range("I15").Value = 0.07155529225243931
$ws.Range("J15").Value = 0.07155529225243933
$ws.Range("O15").Value = 0.1196497582104962
$ws.Range("P15").Value = 0.1196497582104962
$ws.Range("Q15").Value = 0.0508626381591111
$ws.Range("R15").Value = 0.4577637434319999
$ws.Range("S15").Value = 0.008561573416685753
$ws.Range("T15").Value = 0.008561573416685756

$ws.Range("G16").Value = 0.05954566666666666
$ws.Range("H16").Value = 0.178637
$ws.Range("I16").Value = 0.07155529225243931
$ws.Range("J16").Value = 0.07155529225243933
$ws.Range("M16").Value = 0.1824346666666667
$ws.Range("N16").Value = 0.547304
$ws.Range("O16").Value = 0.02555468148257719
$ws.Range("P16").Value = 0.02555468148257719
$ws.Range("Q16").Value = 0.01086319384977778
$ws.Range("R16").Value = 0.09776874464799999
$ws.Range("S16").Value = 0.00182857270190381
$ws.Range("T16").Value = 0.00182857270190381

$ws.Range("G17").Value = 0.05954566666666666
$ws.Range("H17").Value = 0.178637
$ws.Range("I17").Value = 0.07155529225243931
$ws.Range("J17").Value = 0.07155529225243933
$ws.Range("M17").Value = 0.035773
$ws.Range("N17").Value = 0.107319
$ws.Range("O17").Value = 0.005010931515261538
$ws.Range("P17").Value = 0.005010931515261539
$ws.Range("Q17").Value = 0.002130127133666667
$ws.Range("R17").Value = 0.019171144203
$ws.Range("S17").Value = 0.0003585586690314979
$ws.Range("T17").Value = 0.000358558669031498

$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.01206733333333333
$ws.Range("H18").Value = 0.036202
$ws.Range("I18").Value = 0.01450116543673935
$ws.Range("J18").Value = 0.01450116543673936
$ws.Range("M18").Value = 6.066605666666667
$ws.Range("N18").Value = 18.199817
$ws.Range("O18").Value = 0.8497846287916651
$ws.Range("P18").Value = 0.8497846287916652
$ws.Range("Q18").Value = 0.07320775278155556
$ws.Range("R18").Value = 0.6588697750339999
$ws.Range("S18").Value = 0.01232286748770608
$ws.Range("T18").Value = 0.01232286748770608

$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.01206733333333333
$ws.Range("H19").Value = 0.036202
$ws.Range("I19").Value = 0.01450116543673935
$ws.Range("J19").Value = 0.01450116543673936
$ws.Range("O19").Value = 0.1196497582104962
$ws.Range("P19").Value = 0.1196497582104962
$ws.Range("Q19").Value = 0.01030765869688889
$ws.Range("R19").Value = 0.09276892827199998
$ws.Range("S19").Value = 0.001735060938276268
$ws.Range("T19").Value = 0.001735060938276269

$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 0.3333333333333333
$ws.Range("G20").Value = 0.01206733333333333
$ws.Range("H20").Value = 0.036202
$ws.Range("I20").Value = 0.01450116543673935
$ws.Range("J20").Value = 0.01450116543673936
$ws.Range("M20").Value = 0.1824346666666667
$ws.Range("N20").Value = 0.547304
$ws.Range("O20").Value = 0.02555468148257719
$ws.Range("P20").Value = 0.02555468148257719
$ws.Range("Q20").Value = 0.002201499934222222
$ws.Range("R20").Value = 0.019813499408
$ws.Range("S20").Value = 0.0003705726638620315
$ws.Range("T20").Value = 0.0003705726638620316

$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 0.3333333333333333
$ws.Range("G21").Value = 0.01206733333333333
$ws.Range("H21").Value = 0.036202
$ws.Range("I21").Value = 0.01450116543673935
$ws.Range("J21").Value = 0.01450116543673936
$ws.Range("M21").Value = 0.035773
$ws.Range("N21").Value = 0.107319
$ws.Range("O21").Value = 0.005010931515261538
$ws.Range("P21").Value = 0.005010931515261539
$ws.Range("Q21").Value = 0.0004316847153333333
$ws.Range("R21").Value = 0.003885162438
$ws.Range("S21").Value = 0.00007266434689497857
$ws.Range("T21").Value = 0.0000726643468949786
